$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily clear the summary cell (E1) so its shared-string slot is
# freed. This lets the two newly-typed names below reuse that freed
# slot (inserted right after "张少永" and before the summary string),
# matching the shared-string ordering produced by the original edit.
$ws.Range("E1").Value = ""

# Two new keyholders: LuJingyu and DuAngang.
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "卢婧宇"
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "杜昂昂"

# Restore the summary text with the updated counts (used 12 -> 14,
# idle 7 -> 5; total stays 19).
$e1 = $ws.Range("E1")
$e1.Value = "（目前共19把，使用14把，闲置5把）"

# Re-apply the original per-segment fonts so the string keeps its rich
# text run structure instead of collapsing into one plain run.
$e1.Characters(5,2).Font.Name = "Arial"
$e1.Characters(5,2).Font.Size = 10
$e1.Characters(7,4).Font.Name = "Noto Sans CJK SC Regular"
$e1.Characters(7,4).Font.Size = 10
$e1.Characters(11,2).Font.Name = "Arial"
$e1.Characters(11,2).Font.Size = 10
$e1.Characters(13,7).Font.Name = "Noto Sans CJK SC Regular"
$e1.Characters(13,7).Font.Size = 10

# Update the active selection to reflect the last-edited row.
$ws.Range("B16").Select() | Out-Null
